# Fruta / hortaliza, semanal
#
# A new weekly price record was inserted into the "Durazno" (peach) price
# list for Vega Monumental Concepción. This shifts every existing record
# from row 24 onward down by one row (old row 128 -> new row 129), and the
# newly opened row 24 is populated with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 24, pushing rows 24-128
# down to 25-129 (dimension grows from A1:T128 to A1:T129).
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new Kurakata / Especial record.
$ws.Cells.Item(24, 1).Value = 11
$ws.Cells.Item(24, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(24, 3).Value = "Bíobío"
$ws.Cells.Item(24, 4).Value = 44565
$ws.Cells.Item(24, 5).Value = 8
$ws.Cells.Item(24, 6).Value = "Fruta"
$ws.Cells.Item(24, 7).Value = 100103
$ws.Cells.Item(24, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(24, 9).Value = 100103004
$ws.Cells.Item(24, 10).Value = "Durazno"
$ws.Cells.Item(24, 11).Value = "Kurakata"
$ws.Cells.Item(24, 12).Value = "Especial"
$ws.Cells.Item(24, 13).Value = 150
$ws.Cells.Item(24, 14).Value = 14000
$ws.Cells.Item(24, 15).Value = 15000
$ws.Cells.Item(24, 16).Value = 14467
$ws.Cells.Item(24, 17).Value = "$/caja 16 kilos empedrada"
$ws.Cells.Item(24, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(24, 19).Value = 904
$ws.Cells.Item(24, 20).Value = 16
